# pension_zero.xlsx update
#
# 1. Header "Betalt (inkl kostnader)" (col T, row 5) on both the "private"
#    and "Folketrygden" sheets is renamed to "Finans kostnader)".
# 2. On the "Folketrygden" sheet, the AG-column asset-rule labels change:
#       "000 Asset rule "                      -> " Asset rule: "
#       "Folketrygden fra 204200 Asset rule "  -> "Folketrygden fra 2042 Asset rule: "
#    (row 57 - the "2042" row - keeps the longer label, every other data
#    row 6-72 gets the short label).
# 3. Column T (20) on both sheets is narrowed to match column S's width,
#    since the new, shorter header text no longer needs the old width.
# 4. The helper column W (rows 6-72) on "Folketrygden" holds the literal
#    value zero.

$wb = $excel.ActiveWorkbook

$wsPrivate = $wb.Worksheets.Item("private")
$wsPension = $wb.Worksheets.Item("Folketrygden")

# --- 1. Column header text -------------------------------------------------
$wsPrivate.Range("T5").Value = "Finans kostnader)"
$wsPension.Range("T5").Value = "Finans kostnader)"

# --- 2. AG-column asset-rule labels (Folketrygden sheet only) --------------
for ($row = 6; $row -le 72; $row++) {
    if ($row -eq 57) {
        $wsPension.Range("AG57").Value = "Folketrygden fra 2042 Asset rule: "
    } else {
        $wsPension.Range("AG" + $row).Value = " Asset rule: "
    }
}

# --- 3. Column T width, matches column S on both sheets --------------------
$wsPrivate.Columns.Item(20).ColumnWidth = $wsPrivate.Columns.Item(19).ColumnWidth
$wsPension.Columns.Item(20).ColumnWidth = $wsPension.Columns.Item(19).ColumnWidth

# --- 4. W6:W72 literal zero on the Folketrygden sheet -----------------------
for ($row = 6; $row -le 72; $row++) {
    $wsPension.Range("W" + $row).Value = 0
}
